$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$firstRow = 2
$lastRow = 483

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
